$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 163.1857856666667
$ws.Range("H2").Value = 489.557357
$ws.Range("I2").Value = 0.3160920635566714
$ws.Range("J2").Value = 0.3160920635566714
$ws.Range("M2").Value = 70.46291600000001
$ws.Range("N2").Value = 211.388748
$ws.Range("O2").Value = 0.5276750397950939
$ws.Range("P2").Value = 0.5276750397950939
$ws.Range("Q2").Value = 11498.54630782434
$ws.Range("R2").Value = 103486.916770419
$ws.Range("S2").Value = 0.1667938922161799
$ws.Range("T2").Value = 0.1667938922161799
$ws.Range("G3").Value = 163.1857856666667
$ws.Range("H3").Value = 489.557357
$ws.Range("I3").Value = 0.3160920635566714
$ws.Range("J3").Value = 0.3160920635566714
$ws.Range("O3").Value = 0.07361176802536967
$ws.Range("P3").Value = 0.07361176802536967
$ws.Range("Q3").Value = 1604.071179431227
$ws.Range("R3").Value = 14436.64061488105
$ws.Range("S3").Value = 0.0232680956571941
$ws.Range("T3").Value = 0.0232680956571941
$ws.Range("G4").Value = 163.1857856666667
$ws.Range("H4").Value = 489.557357
$ws.Range("I4").Value = 0.3160920635566714
$ws.Range("J4").Value = 0.3160920635566714
$ws.Range("M4").Value = 42.505498
$ws.Range("N4").Value = 127.516494
$ws.Range("O4").Value = 0.3183105613832428
$ws.Range("P4").Value = 0.3183105613832428
$ws.Range("Q4").Value = 6936.293086282928
$ws.Range("R4").Value = 62426.63777654635
$ws.Range("S4").Value = 0.1006154421995117
$ws.Range("T4").Value = 0.1006154421995117
$ws.Range("G5").Value = 163.1857856666667
$ws.Range("H5").Value = 489.557357
$ws.Range("I5").Value = 0.3160920635566714
$ws.Range("J5").Value = 0.3160920635566714
$ws.Range("M5").Value = 10.73653933333333
$ws.Range("N5").Value = 32.209618
$ws.Range("O5").Value = 0.08040263079629371
$ws.Range("P5").Value = 0.08040263079629371
$ws.Range("Q5").Value = 1752.05060645107
$ws.Range("R5").Value = 15768.45545805963
$ws.Range("S5").Value = 0.02541463348378566
$ws.Range("T5").Value = 0.02541463348378566
$ws.Range("I6").Value = 0.1862036026373569
$ws.Range("J6").Value = 0.1862036026373569
$ws.Range("M6").Value = 70.46291600000001
$ws.Range("N6").Value = 211.388748
$ws.Range("O6").Value = 0.5276750397950939
$ws.Range("P6").Value = 0.5276750397950939
$ws.Range("Q6").Value = 6773.566927046565
$ws.Range("R6").Value = 60962.10234341909
$ws.Range("S6").Value = 0.09825499343165714
$ws.Range("T6").Value = 0.09825499343165713
$ws.Range("I7").Value = 0.1862036026373569
$ws.Range("J7").Value = 0.1862036026373569
$ws.Range("O7").Value = 0.07361176802536967
$ws.Range("P7").Value = 0.07361176802536967
$ws.Range("S7").Value = 0.01370677640282923
$ws.Range("T7").Value = 0.01370677640282923
$ws.Range("I8").Value = 0.1862036026373569
$ws.Range("J8").Value = 0.1862036026373569
$ws.Range("M8").Value = 42.505498
$ws.Range("N8").Value = 127.516494
$ws.Range("O8").Value = 0.3183105613832428
$ws.Range("P8").Value = 0.3183105613832428
$ws.Range("Q8").Value = 4086.033502650442
$ws.Range("R8").Value = 36774.30152385398
$ws.Range("S8").Value = 0.05927057328707935
$ws.Range("T8").Value = 0.05927057328707933
$ws.Range("I9").Value = 0.1862036026373569
$ws.Range("J9").Value = 0.1862036026373569
$ws.Range("M9").Value = 10.73653933333333
$ws.Range("N9").Value = 32.209618
$ws.Range("O9").Value = 0.08040263079629371
$ws.Range("P9").Value = 0.08040263079629371
$ws.Range("Q9").Value = 1032.098469203307
$ws.Range("R9").Value = 9288.886222829766
$ws.Range("S9").Value = 0.01497125951579119
$ws.Range("T9").Value = 0.01497125951579119
$ws.Range("G10").Value = 133.6085763333333
$ws.Range("H10").Value = 400.825729
$ws.Range("I10").Value = 0.2588007921740151
$ws.Range("J10").Value = 0.2588007921740151
$ws.Range("M10").Value = 70.46291600000001
$ws.Range("N10").Value = 211.388748
$ws.Range("O10").Value = 0.5276750397950939
$ws.Range("P10").Value = 0.5276750397950939
$ws.Range("Q10").Value = 9414.449891055256
$ws.Range("R10").Value = 84730.04901949731
$ws.Range("S10").Value = 0.1365627183094253
$ws.Range("T10").Value = 0.1365627183094253
$ws.Range("G11").Value = 133.6085763333333
$ws.Range("H11").Value = 400.825729
$ws.Range("I11").Value = 0.2588007921740151
$ws.Range("J11").Value = 0.2588007921740151
$ws.Range("O11").Value = 0.07361176802536967
$ws.Range("P11").Value = 0.07361176802536967
$ws.Range("Q11").Value = 1313.335384853407
$ws.Range("R11").Value = 11820.01846368066
$ws.Range("S11").Value = 0.01905078387829551
$ws.Range("T11").Value = 0.01905078387829551
$ws.Range("G12").Value = 133.6085763333333
$ws.Range("H12").Value = 400.825729
$ws.Range("I12").Value = 0.2588007921740151
$ws.Range("J12").Value = 0.2588007921740151
$ws.Range("M12").Value = 42.505498
$ws.Range("N12").Value = 127.516494
$ws.Range("O12").Value = 0.3183105613832428
$ws.Range("P12").Value = 0.3183105613832428
$ws.Range("Q12").Value = 5679.099074119346
$ws.Range("R12").Value = 51111.89166707412
$ws.Range("S12").Value = 0.08237902544333871
$ws.Range("T12").Value = 0.08237902544333871
$ws.Range("G13").Value = 133.6085763333333
$ws.Range("H13").Value = 400.825729
$ws.Range("I13").Value = 0.2588007921740151
$ws.Range("J13").Value = 0.2588007921740151
$ws.Range("M13").Value = 10.73653933333333
$ws.Range("N13").Value = 32.209618
$ws.Range("O13").Value = 0.08040263079629371
$ws.Range("P13").Value = 0.08040263079629371
$ws.Range("Q13").Value = 1434.493735073502
$ws.Range("R13").Value = 12910.44361566152
$ws.Range("S13").Value = 0.02080826454295568
$ws.Range("T13").Value = 0.02080826454295568
$ws.Range("G14").Value = 123.3364156666667
$ws.Range("H14").Value = 370.009247
$ws.Range("I14").Value = 0.2389035416319566
$ws.Range("J14").Value = 0.2389035416319566
$ws.Range("M14").Value = 70.46291600000001
$ws.Range("N14").Value = 211.388748
$ws.Range("O14").Value = 0.5276750397950939
$ws.Range("P14").Value = 0.5276750397950939
$ws.Range("Q14").Value = 8690.643496861418
$ws.Range("R14").Value = 78215.79147175276
$ws.Range("S14").Value = 0.1260634358378316
$ws.Range("T14").Value = 0.1260634358378316
$ws.Range("G15").Value = 123.3364156666667
$ws.Range("H15").Value = 370.009247
$ws.Range("I15").Value = 0.2389035416319566
$ws.Range("J15").Value = 0.2389035416319566
$ws.Range("O15").Value = 0.07361176802536967
$ws.Range("P15").Value = 0.07361176802536967
$ws.Range("Q15").Value = 1212.362884040471
$ws.Range("R15").Value = 10911.26595636424
$ws.Range("S15").Value = 0.01758611208705083
$ws.Range("T15").Value = 0.01758611208705083
$ws.Range("G16").Value = 123.3364156666667
$ws.Range("H16").Value = 370.009247
$ws.Range("I16").Value = 0.2389035416319566
$ws.Range("J16").Value = 0.2389035416319566
$ws.Range("M16").Value = 42.505498
$ws.Range("N16").Value = 127.516494
$ws.Range("O16").Value = 0.3183105613832428
$ws.Range("P16").Value = 0.3183105613832428
$ws.Range("Q16").Value = 5242.475769446668
$ws.Range("R16").Value = 47182.28192502002
$ws.Range("S16").Value = 0.07604552045331302
$ws.Range("T16").Value = 0.07604552045331302
$ws.Range("G17").Value = 123.3364156666667
$ws.Range("H17").Value = 370.009247
$ws.Range("I17").Value = 0.2389035416319566
$ws.Range("J17").Value = 0.2389035416319566
$ws.Range("M17").Value = 10.73653933333333
$ws.Range("N17").Value = 32.209618
$ws.Range("O17").Value = 0.08040263079629371
$ws.Range("P17").Value = 0.08040263079629371
$ws.Range("Q17").Value = 1324.206278037516
$ws.Range("R17").Value = 11917.85650233765
$ws.Range("S17").Value = 0.01920847325376119
$ws.Range("T17").Value = 0.01920847325376119
